$wb = $excel.ActiveWorkbook

# --- Rename header labels on existing sheets ---
$ws1 = $wb.Worksheets.Item("Weekly Quantity")
$ws1.Range("B1").Value = "Weekly_PO_Qty"

$ws2 = $wb.Worksheets.Item("Monthly Trend")
$ws2.Range("B1").Value = "Monthly_PO_Qty"

# --- Add the new "PO Forecast" sheet after "Monthly Trend" ---
$newSheet = $wb.Worksheets.Add($null, $ws2)
$newSheet.Name = "PO Forecast"

# --- Header row ---
$newSheet.Range("A1").Value = "ds"
$newSheet.Range("B1").Value = "PO_Forecast"
$newSheet.Range("C1").Value = "yhat_lower"
$newSheet.Range("D1").Value = "yhat_upper"

# --- Data rows ---
$newSheet.Range("A2").Value = 45354.99999999999
$newSheet.Range("B2").Value = 6
$newSheet.Range("C2").Value = -20.69809425462256
$newSheet.Range("D2").Value = 32.66900038382212
$newSheet.Range("A3").Value = 45361.99999999999
$newSheet.Range("B3").Value = 6
$newSheet.Range("C3").Value = -18.71328976400946
$newSheet.Range("D3").Value = 34.03058759158448
$newSheet.Range("A4").Value = 45375.99999999999
$newSheet.Range("B4").Value = 7
$newSheet.Range("C4").Value = -18.65519148826692
$newSheet.Range("D4").Value = 33.00894182527389
$newSheet.Range("A5").Value = 45389.99999999999
$newSheet.Range("B5").Value = 9
$newSheet.Range("C5").Value = -16.19287692619292
$newSheet.Range("D5").Value = 35.38629251703927
$newSheet.Range("A6").Value = 45396.99999999999
$newSheet.Range("B6").Value = 9
$newSheet.Range("C6").Value = -17.70985541773679
$newSheet.Range("D6").Value = 35.29944309351595
$newSheet.Range("A7").Value = 45417.99999999999
$newSheet.Range("B7").Value = 11
$newSheet.Range("C7").Value = -12.2329218010709
$newSheet.Range("D7").Value = 38.09091248956442
$newSheet.Range("A8").Value = 45424.99999999999
$newSheet.Range("B8").Value = 12
$newSheet.Range("C8").Value = -15.49974378191732
$newSheet.Range("D8").Value = 36.56458672716339
$newSheet.Range("A9").Value = 45431.99999999999
$newSheet.Range("B9").Value = 12
$newSheet.Range("C9").Value = -12.62939106590607
$newSheet.Range("D9").Value = 38.97549038173011
$newSheet.Range("A10").Value = 45438.99999999999
$newSheet.Range("B10").Value = 13
$newSheet.Range("C10").Value = -12.68394994105923
$newSheet.Range("D10").Value = 38.53133991632953
$newSheet.Range("A11").Value = 45445.99999999999
$newSheet.Range("B11").Value = 14
$newSheet.Range("C11").Value = -11.79482354803515
$newSheet.Range("D11").Value = 38.0977300123448
$newSheet.Range("A12").Value = 45452.99999999999
$newSheet.Range("B12").Value = 14
$newSheet.Range("C12").Value = -10.12460743566958
$newSheet.Range("D12").Value = 40.16068951033957
$newSheet.Range("A13").Value = 45466.99999999999
$newSheet.Range("B13").Value = 16
$newSheet.Range("C13").Value = -9.83485232446562
$newSheet.Range("D13").Value = 42.74667929590973
$newSheet.Range("A14").Value = 45473.99999999999
$newSheet.Range("B14").Value = 16
$newSheet.Range("C14").Value = -8.290469450244728
$newSheet.Range("D14").Value = 44.49613360138262
$newSheet.Range("A15").Value = 45487.99999999999
$newSheet.Range("B15").Value = 17
$newSheet.Range("C15").Value = -7.05472334712988
$newSheet.Range("D15").Value = 44.09416768714909
$newSheet.Range("A16").Value = 45494.99999999999
$newSheet.Range("B16").Value = 18
$newSheet.Range("C16").Value = -9.4979813974264
$newSheet.Range("D16").Value = 41.53676493744047
$newSheet.Range("A17").Value = 45501.99999999999
$newSheet.Range("B17").Value = 19
$newSheet.Range("C17").Value = -6.214833353856172
$newSheet.Range("D17").Value = 43.44465330454558
$newSheet.Range("A18").Value = 45508.99999999999
$newSheet.Range("B18").Value = 19
$newSheet.Range("C18").Value = -7.718148553855489
$newSheet.Range("D18").Value = 42.81040198142225
$newSheet.Range("A19").Value = 45515.99999999999
$newSheet.Range("B19").Value = 20
$newSheet.Range("C19").Value = -5.941160342244919
$newSheet.Range("D19").Value = 46.14373275515684
$newSheet.Range("A20").Value = 45543.99999999999
$newSheet.Range("B20").Value = 22
$newSheet.Range("C20").Value = -1.730924932663882
$newSheet.Range("D20").Value = 47.93782354184591
$newSheet.Range("A21").Value = 45578.99999999999
$newSheet.Range("B21").Value = 26
$newSheet.Range("C21").Value = -0.3173594224663284
$newSheet.Range("D21").Value = 49.32858980400173
$newSheet.Range("A22").Value = 45585.99999999999
$newSheet.Range("B22").Value = 26
$newSheet.Range("C22").Value = -0.3505000201405808
$newSheet.Range("D22").Value = 52.08264221853875
$newSheet.Range("A23").Value = 45606.99999999999
$newSheet.Range("B23").Value = 28
$newSheet.Range("C23").Value = 1.774829029442369
$newSheet.Range("D23").Value = 55.37929238564558
$newSheet.Range("A24").Value = 45634.99999999999
$newSheet.Range("B24").Value = 30
$newSheet.Range("C24").Value = 5.84372305284962
$newSheet.Range("D24").Value = 55.84413267084251
$newSheet.Range("A25").Value = 45641.99999999999
$newSheet.Range("B25").Value = 31
$newSheet.Range("C25").Value = 6.969054057848395
$newSheet.Range("D25").Value = 54.6794552907532
$newSheet.Range("A26").Value = 45648.99999999999
$newSheet.Range("B26").Value = 32
$newSheet.Range("C26").Value = 5.665284879289796
$newSheet.Range("D26").Value = 58.82771623704178
$newSheet.Range("A27").Value = 45655.99999999999
$newSheet.Range("B27").Value = 32
$newSheet.Range("C27").Value = 5.678460292573053
$newSheet.Range("D27").Value = 57.91531124365638
$newSheet.Range("A28").Value = 45662.99999999999
$newSheet.Range("B28").Value = 33
$newSheet.Range("C28").Value = 8.793366434890293
$newSheet.Range("D28").Value = 58.98724103723869
$newSheet.Range("A29").Value = 45669.99999999999
$newSheet.Range("B29").Value = 34
$newSheet.Range("C29").Value = 8.452096523622505
$newSheet.Range("D29").Value = 59.63275629215116
$newSheet.Range("A30").Value = 45676.99999999999
$newSheet.Range("B30").Value = 34
$newSheet.Range("C30").Value = 9.324734872345893
$newSheet.Range("D30").Value = 59.17304874452327
$newSheet.Range("A31").Value = 45683.99999999999
$newSheet.Range("B31").Value = 35
$newSheet.Range("C31").Value = 9.461149311870404
$newSheet.Range("D31").Value = 59.78022319504099
$newSheet.Range("A32").Value = 45690.99999999999
$newSheet.Range("B32").Value = 35
$newSheet.Range("C32").Value = 11.18127552724051
$newSheet.Range("D32").Value = 60.96498032340661

# --- Formatting: reuse the existing header/date styles from "Weekly Quantity" ---
$ws1.Range("A1:B1").Copy()
$newSheet.Range("A1:D1").PasteSpecial(-4122)  # xlPasteFormats

$ws1.Range("A2").Copy()
$newSheet.Range("A2:A32").PasteSpecial(-4122)  # xlPasteFormats

Write-Output "done"
